$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clientes")

$ws.Cells.Item(12, 1).Value = 1000535410
$ws.Cells.Item(12, 2).Value = "sdafsafdsdf"
$ws.Cells.Item(12, 3).Value = 6648792520

$ws.Cells.Item(13, 1).Value = 1000535410
$ws.Cells.Item(13, 2).Value = "aasd"
$ws.Cells.Item(13, 3).Value = 2269735498
